# Apply header restructuring changes described in the commit:
# "fixing sample excel files + report values"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primary Review Worksheet")

# Column M: "Department" -> "Agency / Office"
$ws.Range("M1").Value = "Agency / Office"

# Column U: "Total Funding" -> "Department"
$ws.Range("U1").Value = "Department"

# Columns V..AA: replace old BY funding / Has Keywords headers with new RAI headers
$ws.Range("V1").Value = "RAI Secondary Reviewer"
$ws.Range("W1").Value = "RAI Tag Agree"
$ws.Range("X1").Value = "RAI Tag"
$ws.Range("Y1").Value = "RAI Transition Partner Agree"
$ws.Range("Z1").Value = "RAI Transition Partner"
$ws.Range("AA1").Value = "RAI Mission Partners"

# New columns AB..AG
$ws.Range("AB1").Value = "POC Title"
$ws.Range("AC1").Value = "POC Name"
$ws.Range("AD1").Value = "POC Email"
$ws.Range("AE1").Value = "POC Org"
$ws.Range("AF1").Value = "POC Phone Number"
$ws.Range("AG1").Value = "RAI Review Notes"

# Make sure all the header cells keep/get the same bold/underline header format (style index 1)
# used by the rest of row 1 -- copy the format from A1 onto the touched/extended cells.
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("U1:AG1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view state to match the saved selection/scroll position
$ws.Range("X18").Select()
$ws.Application.ActiveWindow.ScrollColumn = 12
